# "Removing excel.xlsx mapping for custom functions snippets because the
# process doesn't work"
#
# The "Snippets" table on the active sheet has two rows describing the
# excel-custom-functions-errors snippet (Class "Error" / member "class",
# and Class "ErrorCode" / member "excel-custom-functions-errors") at
# worksheet rows 142 and 145. Delete those two rows outright so the table
# (and everything below them) shifts up and the sheet shrinks from
# A1:E370 to A1:E368.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the lower row first so row 142's address is still valid when we
# get to it.
$ws.Rows.Item(145).Delete()
$ws.Rows.Item(142).Delete()

# Match the resulting view/selection state (scrolled down a bit less than
# before, with the whole data row selected where the deletion landed).
$win = $ws.Application.ActiveWindow
$win.ScrollRow = 100
$ws.Range("A118:XFD118").Select()
